$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values
$ws.Range("C11").Value = 8
$ws.Range("C24").Value = $null
$ws.Range("C25").Value = 10
$ws.Range("C27").Value = 10
$ws.Range("C36").Value = 8
$ws.Range("C38").Value = 8
$ws.Range("C39").Value = 5

# Set column F width (new col min=6 max=6 width 8.88671875 characters).
# The host's ColumnWidth setter snaps to 1/6-character increments, so 8 is
# the closest representable value to 8.88671875 (rounds to 8.8333...).
$ws.Columns.Item(6).ColumnWidth = 8

# Update the view selection (was C20, now J4); the saved view also stops
# scrolling to a fixed topLeftCell once the selection is written.
$ws.Range("J4").Select()
